# "Generate Report for Handoff" — refresh the localization-status report
# with the results of the latest handoff/handback run for the two files
# that moved from "Handed back: in sync with en-US" to "Ready for
# handoff", and flag them with a stale-handback-version error detail.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$errorDetail_9b1146d9 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ed8052d97f3c112d855b56cfb59d8c5332db6f4/e2e/9b1146d9-1709-4171-b1ed-39560ae634a7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77e5f5fc3e435cda0eca8d8fc567960ac7f137dd/e2e/9b1146d9-1709-4171-b1ed-39560ae634a7.md."
$errorDetail_9bd58b82 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ed8052d97f3c112d855b56cfb59d8c5332db6f4/e2e/9bd58b82-ca50-42b9-b7dc-6bd02cd37065.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77e5f5fc3e435cda0eca8d8fc567960ac7f137dd/e2e/9bd58b82-ca50-42b9-b7dc-6bd02cd37065.md."

# ---------------------------------------------------------------
# Overview sheet — rows 4 (9b1146d9...) and 5 (9bd58b82...)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus
$wsOverview.Range("G4").Value = "2016-09-05 00:30:31"

$wsOverview.Range("E5").Value = $newStatus
$wsOverview.Range("F5").Value = $newStatus
$wsOverview.Range("G5").Value = "2016-09-05 00:30:31"

# ---------------------------------------------------------------
# zh-cn sheet — rows 4 (9b1146d9...) and 5 (9bd58b82...)
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $newStatus
$wsZhCn.Range("H4").Value = "2016-09-05 00:30:26"
$wsZhCn.Range("P4").Value = $errorDetail_9b1146d9

$wsZhCn.Range("C5").Value = $newStatus
$wsZhCn.Range("H5").Value = "2016-09-05 00:30:26"
$wsZhCn.Range("P5").Value = $errorDetail_9bd58b82

$wsZhCn.Columns.Item(16).ColumnWidth = 39.2

# ---------------------------------------------------------------
# de-de sheet — rows 4 (9b1146d9...) and 5 (9bd58b82...)
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $newStatus
$wsDeDe.Range("H4").Value = "2016-09-05 00:30:31"
$wsDeDe.Range("P4").Value = $errorDetail_9b1146d9

$wsDeDe.Range("C5").Value = $newStatus
$wsDeDe.Range("H5").Value = "2016-09-05 00:30:31"
$wsDeDe.Range("P5").Value = $errorDetail_9bd58b82

$wsDeDe.Columns.Item(16).ColumnWidth = 39.2
